$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("H19").Value = 14706562
$ws.Range("I19").Value = 50000400
$ws.Range("J19").Value = 796.75
$ws.Range("K19").Value = 50000400
$ws.Range("L19").Value = 796.75
$ws.Range("M19").Value = -50000225
$ws.Range("N19").Value = -1146.75
$ws.Range("H41").Value = 11765033
$ws.Range("I41").Value = 458.875
$ws.Range("J41").Value = 22222432
$ws.Range("K41").Value = 458.875
$ws.Range("L41").Value = 22222432
$ws.Range("M41").Value = -18.875
$ws.Range("N41").Value = -22223312
$ws.Range("H86").Value = 27125.75
$ws.Range("I86").Value = 34967.668
$ws.Range("J86").Value = 3600
$ws.Range("K86").Value = 34967.668
$ws.Range("L86").Value = 3600
$ws.Range("M86").Value = -33844.668
$ws.Range("N86").Value = -5846
$ws.Range("H89").Value = 27125.75
$ws.Range("I89").Value = 34967.668
$ws.Range("J89").Value = 3600
$ws.Range("K89").Value = 174838.34
$ws.Range("L89").Value = 18000
$ws.Range("M89").Value = -169222.34
$ws.Range("N89").Value = -29232
$ws.Range("H108").Value = 34842
$ws.Range("J108").Value = 34842
$ws.Range("L108").Value = 34842
$ws.Range("N108").Value = -42522
$ws.Range("H132").Value = 7147611.5
$ws.Range("I132").Value = 8004849
$ws.Range("K132").Value = 24014547
$ws.Range("M132").Value = -24012017
$ws.Range("H135").Value = 716.4400000000001
$ws.Range("I135").Value = 765.04346
$ws.Range("K135").Value = 6885.39114
$ws.Range("M135").Value = -4350.39114

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2943.9167
$ws.Range("I61").Value = 1578.8235
$ws.Range("K61").Value = 1578.8235
$ws.Range("M61").Value = -1366.8235
$ws.Range("H75").Value = 38000
$ws.Range("J75").Value = 38000
$ws.Range("L75").Value = 38000
$ws.Range("N75").Value = -39748
$ws.Range("H78").Value = 38000
$ws.Range("J78").Value = 38000
$ws.Range("L78").Value = 114000
$ws.Range("N78").Value = -122736
$ws.Range("H97").Value = 524.75
$ws.Range("I97").Value = 482
$ws.Range("K97").Value = 482
$ws.Range("M97").Value = 14
$ws.Range("H98").Value = 29999.875
$ws.Range("J98").Value = 29999.875
$ws.Range("L98").Value = 29999.875
$ws.Range("N98").Value = -35989.875
$ws.Range("H110").Value = 1122.0588
$ws.Range("I110").Value = 540.3226
$ws.Range("K110").Value = 540.3226
$ws.Range("M110").Value = 1504.6774
$ws.Range("H136").Value = 2943.9167
$ws.Range("I136").Value = 1578.8235
$ws.Range("K136").Value = 4736.470499999999
$ws.Range("M136").Value = -2186.470499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16131768
$ws.Range("I58").Value = 1821
$ws.Range("J58").Value = 71434440
$ws.Range("K58").Value = 1821
$ws.Range("L58").Value = 71434440
$ws.Range("M58").Value = -1618
$ws.Range("N58").Value = -71434846
$ws.Range("H105").Value = 2359.2144
$ws.Range("I105").Value = 1502.9
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 1502.9
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = 244.0999999999999
$ws.Range("N105").Value = -7994
$ws.Range("H132").Value = 2598.1167
$ws.Range("I132").Value = 1752.2222
$ws.Range("J132").Value = 3866.9583
$ws.Range("K132").Value = 5256.6666
$ws.Range("L132").Value = 11600.8749
$ws.Range("M132").Value = -2726.6666
$ws.Range("N132").Value = -16660.8749
$ws.Range("H136").Value = 16131768
$ws.Range("I136").Value = 1821
$ws.Range("J136").Value = 71434440
$ws.Range("K136").Value = 5463
$ws.Range("L136").Value = 214303320
$ws.Range("M136").Value = -2913
$ws.Range("N136").Value = -214308420

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 170.8421
$ws.Range("I12").Value = 16.2
$ws.Range("J12").Value = 226.07143
$ws.Range("K12").Value = 48.59999999999999
$ws.Range("L12").Value = 678.21429
$ws.Range("M12").Value = 124.4
$ws.Range("N12").Value = -1024.21429
$ws.Range("H86").Value = 1500
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -1814
$ws.Range("N86").Value = -8372
$ws.Range("H89").Value = 1500
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 18000
$ws.Range("M89").Value = -3072
$ws.Range("N89").Value = -29856
$ws.Range("H97").Value = 513.4286
$ws.Range("I97").Value = 300
$ws.Range("J97").Value = 549
$ws.Range("K97").Value = 900
$ws.Range("L97").Value = 1647
$ws.Range("M97").Value = -404
$ws.Range("N97").Value = -2639
$ws.Range("H98").Value = 233.33333
$ws.Range("I98").Value = 233.33333
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 699.99999
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 798.00001
$ws.Range("N98").ClearContents()
$ws.Range("H105").Value = 4114.2856
$ws.Range("I105").Value = 4900
$ws.Range("J105").Value = 3983.3333
$ws.Range("K105").Value = 14700
$ws.Range("L105").Value = 11949.9999
$ws.Range("M105").Value = -12079
$ws.Range("N105").Value = -17191.9999
$ws.Range("H107").Value = 1541.4
$ws.Range("J107").Value = 2332.3333
$ws.Range("L107").Value = 6996.999899999999
$ws.Range("N107").Value = -10836.9999
$ws.Range("H110").Value = 3688.0833
$ws.Range("I110").Value = 1625
$ws.Range("J110").Value = 3875.6365
$ws.Range("K110").Value = 4875
$ws.Range("L110").Value = 11626.9095
$ws.Range("M110").Value = -785
$ws.Range("N110").Value = -19806.9095
$ws.Range("H131").Value = 917.0635
$ws.Range("I131").Value = 776.8570999999999
$ws.Range("J131").Value = 1092.3214
$ws.Range("K131").Value = 2330.5713
$ws.Range("L131").Value = 3276.9642
$ws.Range("M131").Value = 2709.4287
$ws.Range("N131").Value = -13356.9642
$ws.Range("H132").Value = 4916.6665
$ws.Range("J132").Value = 5900
$ws.Range("L132").Value = 53100
$ws.Range("N132").Value = -58160
$ws.Range("H133").Value = 3158.8462
$ws.Range("I133").Value = 2649.875
$ws.Range("J133").Value = 3973.2
$ws.Range("K133").Value = 7949.625
$ws.Range("L133").Value = 11919.6
$ws.Range("M133").Value = -2889.625
$ws.Range("N133").Value = -22039.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 12533955
$ws.Range("I24").Value = 5806
$ws.Range("J24").Value = 16710005
$ws.Range("K24").Value = 5806
$ws.Range("L24").Value = 16710005
$ws.Range("M24").Value = -5633
$ws.Range("N24").Value = -16710351
$ws.Range("H70").Value = 3818.3845
$ws.Range("I70").Value = 3761.8
$ws.Range("J70").Value = 3853.75
$ws.Range("K70").Value = 3761.8
$ws.Range("L70").Value = 3853.75
$ws.Range("M70").Value = -3491.8
$ws.Range("N70").Value = -4393.75
$ws.Range("H73").Value = 3818.3845
$ws.Range("I73").Value = 3761.8
$ws.Range("J73").Value = 3853.75
$ws.Range("K73").Value = 3761.8
$ws.Range("L73").Value = 3853.75
$ws.Range("M73").Value = -2825.8
$ws.Range("N73").Value = -5725.75
$ws.Range("H75").Value = 38666.668
$ws.Range("J75").Value = 38666.668
$ws.Range("L75").Value = 38666.668
$ws.Range("N75").Value = -40414.668
$ws.Range("H78").Value = 38666.668
$ws.Range("J78").Value = 38666.668
$ws.Range("L78").Value = 116000.004
$ws.Range("N78").Value = -124736.004
$ws.Range("H122").Value = 7810
$ws.Range("I122").Value = 8350
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 25050
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -22600
$ws.Range("N122").Value = -25900
$ws.Range("H137").Value = 29833.334
$ws.Range("J137").Value = 29833.334
$ws.Range("L137").Value = 29833.334
$ws.Range("N137").Value = -40033.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 324879
$ws.Range("I14").Value = 628503
$ws.Range("J14").Value = 21255
$ws.Range("K14").Value = 628503
$ws.Range("L14").Value = 21255
$ws.Range("M14").Value = -628331
$ws.Range("N14").Value = -21599
$ws.Range("H22").Value = 111112616
$ws.Range("I22").Value = 200000340
$ws.Range("J22").Value = 2952.5
$ws.Range("K22").Value = 200000340
$ws.Range("L22").Value = 2952.5
$ws.Range("M22").Value = -200000045
$ws.Range("N22").Value = -3542.5
$ws.Range("H27").Value = 111112616
$ws.Range("I27").Value = 200000340
$ws.Range("J27").Value = 2952.5
$ws.Range("K27").Value = 200000340
$ws.Range("L27").Value = 2952.5
$ws.Range("M27").Value = -200000233
$ws.Range("N27").Value = -3166.5
$ws.Range("H46").Value = 2417.5
$ws.Range("I46").Value = 826.25
$ws.Range("J46").Value = 3478.3333
$ws.Range("K46").Value = 826.25
$ws.Range("L46").Value = 3478.3333
$ws.Range("M46").Value = -638.25
$ws.Range("N46").Value = -3854.3333
$ws.Range("H122").Value = 3336.3333
$ws.Range("I122").Value = 2592.5557
$ws.Range("J122").Value = 6683.3335
$ws.Range("K122").Value = 7777.6671
$ws.Range("L122").Value = 20050.0005
$ws.Range("M122").Value = -5327.6671
$ws.Range("N122").Value = -24950.0005
$ws.Range("H133").Value = 29661.666
$ws.Range("J133").Value = 29661.666
$ws.Range("L133").Value = 29661.666
$ws.Range("N133").Value = -34721.666
$ws.Range("H136").Value = 3633.3333
$ws.Range("I136").Value = 2900
$ws.Range("J136").Value = 4122.222
$ws.Range("K136").Value = 8700
$ws.Range("L136").Value = 12366.666
$ws.Range("M136").Value = -6150
$ws.Range("N136").Value = -17466.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 57505.25
$ws.Range("J15").Value = 57505.25
$ws.Range("L15").Value = 57505.25
$ws.Range("N15").Value = -58081.25
$ws.Range("H135").Value = 58969.168
$ws.Range("J135").Value = 58969.168
$ws.Range("L135").Value = 58969.168
$ws.Range("N135").Value = -69109.16800000001
$ws.Range("H136").Value = 2507.7856
$ws.Range("I136").Value = 1567.3334
$ws.Range("J136").Value = 3213.125
$ws.Range("K136").Value = 4702.0002
$ws.Range("L136").Value = 9639.375
$ws.Range("M136").Value = -2152.0002
$ws.Range("N136").Value = -14739.375
